$wb = $excel.ActiveWorkbook

# --- Sheet: Scenario_UserSignIn ---
# Reorder merge cells to match target (unmerge + remerge in new order)
$wsScenario = $wb.Worksheets.Item("Scenario_UserSignIn")
$mergeOrder = @("A2:D2","E2:F2","A3:D3","E3:F3","A1:F1")
foreach ($r in $mergeOrder) {
    $wsScenario.Range($r).UnMerge()
}
foreach ($r in $mergeOrder) {
    $wsScenario.Range($r).Merge()
}

# --- Sheet: TC_UserSignIn ---
$ws = $wb.Worksheets.Item("TC_UserSignIn")

# Row height changes
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(53).RowHeight = 15

# Sheet view: scroll position + active selection
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G65").Select()

# Cell value updates (Actual/Expected/Input/etc. text refresh)
$ws.Range('G2').Value = 'City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range('K2').Value = 'title should be "City Super Market Norwalk - Online Grocery Shopping with Home Delivery"'
$ws.Range('L2').Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range('L3').Value = 'VerifyElement: null'
$ws.Range('L4').Value = 'Click: null'
$ws.Range('L5').Value = 'VerifyElement: null'
$ws.Range('L6').Value = 'VerifyText: Sign in'
$ws.Range('L7').Value = 'VerifyElement: null'
$ws.Range('L8').Value = 'VerifyElement: null'
$ws.Range('L9').Value = 'VerifyElement: null'
$ws.Range('L10').Value = 'VerifyElement: null'
$ws.Range('L11').Value = 'VerifyElement: null'
$ws.Range('L12').Value = 'ClearText: null'
$ws.Range('L13').Value = 'ClearText: null'
$ws.Range('L14').Value = 'Click: null'
$ws.Range('L15').Value = 'VerifyText: Please enter your email address.'
$ws.Range('L16').Value = 'VerifyText: Please enter your password.'
$ws.Range('L17').Value = 'ClearText: null'
$ws.Range('L18').Value = 'SetText: 123456'
$ws.Range('L19').Value = 'Click: null'
$ws.Range('L20').Value = 'VerifyText: Please enter your email address.'
$ws.Range('L21').Value = 'VerifyNoElement: null'
$ws.Range('L22').Value = 'SetText: userchandna4861983@mailinator.com'
$ws.Range('L23').Value = 'ClearText: null'
$ws.Range('L24').Value = 'Click: null'
$ws.Range('L25').Value = 'VerifyNoElement: null'
$ws.Range('L26').Value = 'VerifyText: Please enter your password.'
$ws.Range('L27').Value = 'SetText: asdasdasdasd@mailinator.com'
$ws.Range('L28').Value = 'SetText: 123456'
$ws.Range('L29').Value = 'Click: null'
$ws.Range('L30').Value = 'VerifyText: Invalid login or password.'
$ws.Range('L31').Value = 'SetText: userchandna4861983@mailinator.com'
$ws.Range('L32').Value = 'SetText: 12345677'
$ws.Range('L33').Value = 'Click: null'
$ws.Range('L34').Value = 'VerifyText: Invalid login or password.'
$ws.Range('L35').Value = 'SetText: asdasdasdasd@mailinator.com'
$ws.Range('L36').Value = 'SetText: 123456787'
$ws.Range('L37').Value = 'Click: null'
$ws.Range('L38').Value = 'VerifyText: Invalid login or password.'
$ws.Range('L39').Value = 'SetText: userchandna4861983@mailinator.com'
$ws.Range('L40').Value = 'SetText: 123'
$ws.Range('L41').Value = 'Click: null'
$ws.Range('L42').Value = 'VerifyText: Invalid login or password.'
$ws.Range('L43').Value = 'SetText: asdasdasdasd@'
$ws.Range('L44').Value = 'SetText: 123456'
$ws.Range('L45').Value = 'Click: null'
$ws.Range('L46').Value = 'VerifyText: Please enter a valid email address (Ex: johndoe@domain.com).'
$ws.Range('L47').Value = 'SetText: Randomemailid'
$ws.Range('L48').Value = 'SetText: 123456'
$ws.Range('L49').Value = 'Click: null'
$ws.Range('F50').Value = '//span[contains(text(),''Akash'')]'
$ws.Range('L50').Value = 'VerifyText: Akash Sangal'
$ws.Range('L51').Value = 'Click: null'
$ws.Range('L52').Value = 'Click: null'
$ws.Range('L53').Value = 'Wait: 6000'
$ws.Range('L54').Value = 'VerifyElement: null'
$ws.Range('L55').Value = 'Click: null'
$ws.Range('L56').Value = 'Click: null'
$ws.Range('L57').Value = 'VerifyNoElement: null'
$ws.Range('L58').Value = 'Click: null'
$ws.Range('L59').Value = 'Click: null'
$ws.Range('L60').Value = 'VerifyElement: null'
$ws.Range('L61').Value = 'Click: null'
$ws.Range('L62').Value = 'VerifyText: Please enter your email address.'
$ws.Range('L63').Value = 'SetText: Randomemailid'
$ws.Range('L64').Value = 'Click: null'
$ws.Range('C65').Value = 'Verify message "We sent an email with instructions to reset your password. Please check your "Junk" or "Spam" email folders if you do not see the email within the next 10 minutes."'
$ws.Range('G65').Value = 'We sent an email with instructions to reset your password. Please check your "Junk" or "Spam" email folders if you do not see the email within the next 10 minutes'
$ws.Range('K65').Value = 'Message "We sent an email with instructions to reset your password. Please check your "Junk" or "Spam" email folders if you do not see the email within the next 10 minutes." should appear on the screen.'
$ws.Range('L65').Value = 'VerifyText: We sent an email with instructions to reset your password. Please check your Junk or Spam email folders if you do not see the email within the next 10 minutes'
$ws.Range('C66').Value = 'close forgot password form '
$ws.Range('F66').Value = '//div[normalize-space(@id) = ''sociallogin-close-popup'']'
$ws.Range('K66').Value = 'user should able to click '
$ws.Range('L66').Value = 'Click: null'
$ws.Range('B67').Value = 'TC_CityMarket_10'
$ws.Range('L67').Value = 'VerifyNoElement: null'
